# Apply the diff: replace the division expressions in the table cells.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "386÷3="; new = "720÷2="},
    @{old = "451÷8="; new = "679÷7="},
    @{old = "660÷6="; new = "567÷9="},
    @{old = "446÷3="; new = "553÷8="},
    @{old = "133÷6="; new = "184÷3="},
    @{old = "455÷4="; new = "994÷7="},
    @{old = "343÷2="; new = "849÷8="},
    @{old = "392÷4="; new = "940÷7="},
    @{old = "695÷9="; new = "566÷2="},
    @{old = "810÷4="; new = "708÷9="},
    @{old = "303÷7="; new = "143÷9="},
    @{old = "612÷8="; new = "459÷8="},
    @{old = "960÷9="; new = "403÷9="},
    @{old = "743÷5="; new = "797÷8="},
    @{old = "409÷5="; new = "676÷3="},
    @{old = "749÷3="; new = "422÷3="},
    @{old = "943÷8="; new = "183÷6="},
    @{old = "371÷5="; new = "514÷4="},
    @{old = "265÷7="; new = "782÷6="},
    @{old = "519÷2="; new = "988÷2="},
    @{old = "875÷5="; new = "366÷4="},
    @{old = "100÷9="; new = "707÷6="},
    @{old = "356÷5="; new = "445÷5="},
    @{old = "477÷6="; new = "642÷9="},
    @{old = "296÷5="; new = "904÷4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}

$d.Save()
